$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 49
$row = 50

# --- Values ---
$ws.Cells.Item($row, 1).Value = 49
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45224.6875
$ws.Cells.Item($row, 6).Value = "Sumqayit"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Sabah Baku"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 7.39
$ws.Cells.Item($row, 11).Value = "12/08/2023 08:12"
$ws.Cells.Item($row, 12).Value = 6.04
$ws.Cells.Item($row, 13).Value = "25/10/2023 15:39"
$ws.Cells.Item($row, 14).Value = 4.58
$ws.Cells.Item($row, 15).Value = "12/08/2023 08:12"
$ws.Cells.Item($row, 16).Value = 4.06
$ws.Cells.Item($row, 17).Value = "25/10/2023 15:39"
$ws.Cells.Item($row, 18).Value = 1.31
$ws.Cells.Item($row, 19).Value = "12/08/2023 08:12"
$ws.Cells.Item($row, 20).Value = 1.52
$ws.Cells.Item($row, 21).Value = "25/10/2023 15:39"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sumqayit-fk-sabah-baku/n9AWonXS/"

# --- Formatting: copy the formatting from the row above so styles (bold index col, date format) match ---
$ws.Cells.Item($srcRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($srcRow, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = $false
